# Add Area / Atotal columns to the Q discharge sheet (Station 5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ------------------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Row 2: seed area formula + totals side panel -------------------------
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Row 3: explicit (non-shared) area formula -----------------------------
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# --- Rows 4-15: shared area formula (Excel AutoFill adjusts refs per row) -
for ($r = 4; $r -le 15; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 7).Formula = "=(D$r-D$prev)*B$r/100"
}

# Excel re-dimensions the sheet / keeps the selection on the newly-added
# totals cells after doing this kind of edit.
$ws.Range("J2:K2").Select()
